# Hortaliza, Femacal de La Calera - Zanahoria: weekly price update.
# Two new daily price records are inserted into the existing data table:
#   - one after the current row 164 (becomes the new row 165)
#   - one after the current row 259 (becomes row 260 once the sheet has
#     already grown by the first insertion)
# Every row below each insertion point shifts down by one, which is why
# the dimension grows from A1:R267 to A1:R269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

function Set-PriceRow($rowIndex, $values) {
    $ws.Rows.Item($rowIndex).Insert()
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value2 = $values[$col - 1]
    }
    $ws.Cells.Item($rowIndex, 4).NumberFormat = $dateFormat
}

# New record inserted as row 165 (date serial 44567 = 2022-01-06).
Set-PriceRow 165 @(
    3,
    "Femacal de La Calera",
    "Coquimbo",
    44567,
    5,
    100114013,
    "Zanahoria",
    "Sin especificar",
    "Primera",
    400,
    7500,
    8000,
    7775,
    "`$/saco 20 kilos",
    "Provincia de Quillota",
    389,
    20,
    "Hortaliza"
)

# New record inserted as row 260 (date serial 44568 = 2022-01-07), once
# the sheet already reflects the first insertion above.
Set-PriceRow 260 @(
    3,
    "Femacal de La Calera",
    "Coquimbo",
    44568,
    5,
    100114013,
    "Zanahoria",
    "Sin especificar",
    "Primera",
    480,
    7500,
    8000,
    7740,
    "`$/saco 20 kilos",
    "Provincia de Quillota",
    387,
    20,
    "Hortaliza"
)
